$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet:
#    - the existing data row (was "2022-Q2") becomes the new "2022-Q3" row
#    - a fresh row is appended below it with the old "2022-Q2" figures
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q3"

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.01

# Match row 3's index-cell look to row 2's (bold/bordered/centered style)
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Insert a brand-new worksheet for the 2022-Q3 fund breakdown, placed
#    right after "总计" (so the former "2022-Q2" sheet slides down to 3rd).
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q3"

# Mark the fund-code / percentage columns as Text BEFORE writing into them,
# so numeric-looking strings (leading-zero codes, "92.85", …) are kept as
# literal text instead of being coerced to numbers.
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "010343"
$newSheet.Range("C2").Value = "华宝英国富时100指数（QDII）A"
$newSheet.Range("D2").Value = "0.13"
$newSheet.Range("E2").Value = "92.85"
$newSheet.Range("F2").Value = "2.75"
$newSheet.Range("G2").Value = "0.0036"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "010344"
$newSheet.Range("C3").Value = "华宝英国富时100指数（QDII）C"
$newSheet.Range("D3").Value = "0.08"
$newSheet.Range("E3").Value = "92.85"
$newSheet.Range("F3").Value = "2.75"
$newSheet.Range("G3").Value = "0.0022"
$newSheet.Range("H3").Value = 10

# Reuse the same header / index-column style as the "总计" sheet (style
# index 2: bold, thin border, centered) instead of minting a brand-new one.
# This is a format-only paste, applied after the values above, so it cannot
# re-coerce the text cells back into numbers.
$total.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Keep "总计" as the active tab, matching the workbook's original view
#    state (it was untouched by this change).
# ---------------------------------------------------------------------------
$total.Activate()
[void]$total.Range("A1").Select()

Write-Output "2022-Q3 sheet added"
